# #5: property boat&car done
# Fix the "汽車" (car) sheet: row 1 was accidentally populated with the
# same values as the data row instead of real column headers, and the
# data row (row 2) was missing the "capacity" column plus the trailing
# property/category/date/legislator/source/index columns that every
# other property sheet in this workbook carries.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("汽車")

# --- Row 1: proper headers -------------------------------------------------
$ws.Cells.Item(1, 2).Value = "name"
$ws.Cells.Item(1, 3).Value = "capacity"
$ws.Cells.Item(1, 4).Value = "owner"
$ws.Cells.Item(1, 5).Value = "register_date"
$ws.Cells.Item(1, 6).Value = "register_reason"
$ws.Cells.Item(1, 7).Value = "acquire_value"
$ws.Cells.Item(1, 8).Value = "property_category"
$ws.Cells.Item(1, 9).Value = "category"
$ws.Cells.Item(1, 10).Value = "date"
$ws.Cells.Item(1, 11).Value = "legislator_name"
$ws.Cells.Item(1, 12).Value = "legislator_id"
$ws.Cells.Item(1, 13).Value = "source_file"
$ws.Cells.Item(1, 14).Value = "index"

# --- Row 2: the actual car record, now complete -----------------------------
$ws.Cells.Item(2, 1).Value = 47
$ws.Cells.Item(2, 2).Value = "國瑞Camry"
$ws.Cells.Item(2, 3).Value = 1998
$ws.Cells.Item(2, 4).Value = "陳碧涵"
$ws.Cells.Item(2, 5).Value = "96年01月12日"
$ws.Cells.Item(2, 6).Value = "賈賣"
$ws.Cells.Item(2, 7).Value = "(超過五年）"
$ws.Cells.Item(2, 8).Value = "land"
$ws.Cells.Item(2, 9).Value = "normal"
# force this as literal text, not an auto-converted date serial
$ws.Cells.Item(2, 10).NumberFormat = "@"
$ws.Cells.Item(2, 10).Value = "2013-12-31"
$ws.Cells.Item(2, 11).Value = "陳碧涵"
$ws.Cells.Item(2, 12).Value = 1752
$ws.Cells.Item(2, 13).Value = "tmp11ae1"
$ws.Cells.Item(2, 14).Value = 47
